$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended to the CLNERandom sheet (rows 17 and 18)
$rows = @(
    @(42602.010104166664, "Random", 0, 0, 0, 0, 0, 2,  98, 0, 0, 49, 51),
    @(42602.481909722221, "Random", 0, 0, 0, 0, 0, 92, 8,  0, 0, 66, 34)
)

$startRow = 17
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
